# Add 2022-Q4 data
# -----------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计", holding the
#    same shape/formatting as the existing quarter sheets (fund
#    holdings table), filled with the new quarter's data.
# 2) Insert a new summary row into "总计" (shifting the existing rows
#    down) for the new "2022-Q4" quarter.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q4" worksheet right after "总计" -------
# (Grab sheet references AFTER mutating the Worksheets collection --
# references taken before Add()/rename can end up stale.)
$firstSheet = $wb.Worksheets.Item(1)
$newWs = $wb.Worksheets.Add($null, $firstSheet)
$newWs.Name = "2022-Q4"

$summary = $wb.Worksheets.Item("总计")
$q2sheet = $wb.Worksheets.Item("2022-Q2")

# Copy the layout/formatting (header row + first-column style) from
# the "2022-Q2" sheet so the new sheet matches the others.
$q2sheet.Range("A1:H5").Copy()
$newWs.Range("A1").PasteSpecial(-4122)

# Header row (text, keeps the pasted bold/border style s=2)
$newWs.Cells.Item(1, 2).Value2 = "基金代码"
$newWs.Cells.Item(1, 3).Value2 = "基金名称"
$newWs.Cells.Item(1, 4).Value2 = "基金规模"
$newWs.Cells.Item(1, 5).Value2 = "股票总仓位"
$newWs.Cells.Item(1, 6).Value2 = "仓位占比"
$newWs.Cells.Item(1, 7).Value2 = "持有市值(亿元)"
$newWs.Cells.Item(1, 8).Value2 = "仓位排名"

function Set-TextCell($range, $text) {
    # Force text storage (matches the source data, which keeps codes
    # like "004332" / ratios like "0.49" as literal strings rather
    # than numbers), then drop the number-format override so the
    # cell is left with no explicit style -- same as its neighbours.
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.ClearFormats()
}

function Set-FundRow($ws, $row, $idx, $code, $name, $scale, $pos, $ratio, $mv, $rank) {
    $ws.Cells.Item($row, 1).Value2 = $idx
    Set-TextCell $ws.Cells.Item($row, 2) $code
    Set-TextCell $ws.Cells.Item($row, 3) $name
    Set-TextCell $ws.Cells.Item($row, 4) $scale
    Set-TextCell $ws.Cells.Item($row, 5) $pos
    Set-TextCell $ws.Cells.Item($row, 6) $ratio
    Set-TextCell $ws.Cells.Item($row, 7) $mv
    $ws.Cells.Item($row, 8).Value2 = $rank
}

Set-FundRow $newWs 2 0 "004332" "恒生前海沪港深新兴产业精选混合" "0.49" "75.81" "4.80" "0.0235" 4
Set-FundRow $newWs 3 1 "003456" "信澳新目标灵活配置混合"         "0.44" "51.24" "1.30" "0.0057" 4
Set-FundRow $newWs 4 2 "013383" "恒生前海高端制造混合A"           "0.11" "84.98" "3.96" "0.0044" 8
Set-FundRow $newWs 5 3 "013384" "恒生前海高端制造混合C"           "0.04" "84.98" "3.96" "0.0016" 8

# --- 2. Insert the new summary row in "总计" --------------------------
# Shift existing data rows 2..6 down to 3..7 (bottom-up so nothing is
# clobbered before it is read), preserving the first-column style.

for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $summary.Cells.Item($dest, 2).Value2 = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($dest, 3).Value2 = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($dest, 4).Value2 = $summary.Cells.Item($r, 4).Value2
}

# Column A carries the bordered/bold style; copy that format down too.
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)
$summary.Cells.Item(7, 1).Value2 = 5

# New row 2: the 2022-Q4 summary entry
$summary.Cells.Item(2, 2).Value2 = "2022-Q4"
$summary.Cells.Item(2, 3).Value2 = 4
$summary.Cells.Item(2, 4).Value2 = 0.04
